# Robotic Process Automation PMO - add PMO tracking tables
#
# Slide 2 (index 2): add a 2-column "Summary / Assignee" table.
# Slide 3 (index 3): rename the existing risk table, then add three
#                     identical 3-column "Risk Description / Severity / Owner"
#                     risk tables underneath it.
#
# Geometry note: this COM host takes Shapes.AddTable / Left / Top / Width /
# Height in points, and stores EMU = points * 12700 internally. 2032000 EMU
# -> 160.0pt, 719666 EMU -> 56.666614173228346pt, 8128000 EMU -> 640.0pt,
# 5418667 EMU -> 426.6666929133858pt (exact through the AddTable ctor path),
# 370840 EMU -> 29.2pt row height.

$p = $ppt.ActivePresentation

$tableLeft   = 160.0
$tableTop    = 56.666614173228346
$tableWidth  = 640.0
$tableHeight = 426.6666929133858
$rowHeightPt = 29.2

function Add-SummaryTable($slide) {
    $gf = $slide.Shapes.AddTable(3, 2, $tableLeft, $tableTop, $tableWidth, $tableHeight)
    $tbl = $gf.Table

    $tbl.Cell(1, 1).Shape.TextFrame.TextRange.Text = "Summary"
    $tbl.Cell(1, 2).Shape.TextFrame.TextRange.Text = "Assignee"

    $tbl.Cell(2, 1).Shape.TextFrame.TextRange.Text = "Develop Test Case Scenarios"
    $tbl.Cell(2, 2).Shape.TextFrame.TextRange.Text = "elevatebot"

    $tbl.Cell(3, 1).Shape.TextFrame.TextRange.Text = "Create Product Roadmap"
    $tbl.Cell(3, 2).Shape.TextFrame.TextRange.Text = "elevatebot"

    for ($i = 1; $i -le 3; $i++) {
        $tbl.Rows.Item($i).Height = $rowHeightPt
    }

    return $gf
}

function Add-RiskTable($slide) {
    $gf = $slide.Shapes.AddTable(3, 3, $tableLeft, $tableTop, $tableWidth, $tableHeight)
    $tbl = $gf.Table

    $tbl.Cell(1, 1).Shape.TextFrame.TextRange.Text = "Risk Description"
    $tbl.Cell(1, 2).Shape.TextFrame.TextRange.Text = "Severity"
    $tbl.Cell(1, 3).Shape.TextFrame.TextRange.Text = "Owner"

    $tbl.Cell(2, 1).Shape.TextFrame.TextRange.Text = "If something happens, then something bad will happen to the program"
    $tbl.Cell(2, 2).Shape.TextFrame.TextRange.Text = "High"
    $tbl.Cell(2, 3).Shape.TextFrame.TextRange.Text = "Elevate Bot"

    $tbl.Cell(3, 1).Shape.TextFrame.TextRange.Text = "If the bot gets a mind of its own, then we will need to kill it"
    $tbl.Cell(3, 2).Shape.TextFrame.TextRange.Text = "Medium"
    $tbl.Cell(3, 3).Shape.TextFrame.TextRange.Text = "Elevate Bot"

    for ($i = 1; $i -le 3; $i++) {
        $tbl.Rows.Item($i).Height = $rowHeightPt
    }

    return $gf
}

# --- Slide 2: Summary / Assignee table -------------------------------------
$slide2 = $p.Slides.Item(2)
Add-SummaryTable $slide2 | Out-Null

# --- Slide 3: rename existing table, add the three risk tables -------------
$slide3 = $p.Slides.Item(3)

for ($i = 1; $i -le $slide3.Shapes.Count; $i++) {
    $shp = $slide3.Shapes.Item($i)
    if ($shp.HasTable -and $shp.Name -eq "Table 10") {
        $shp.Name = "10"
        break
    }
}

Add-RiskTable $slide3 | Out-Null
Add-RiskTable $slide3 | Out-Null
Add-RiskTable $slide3 | Out-Null
